$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
  'param_TimeStep_starting_index'
  'param_demand1_op_cost_starting_index'
  'param_demand1_inv_cost_starting_index'
  'param_demand2_inv_cost_starting_index'
  'param_demand2_op_cost_starting_index'
  'param_Q_net1_demand2_starting_index'
  'param_net1_sell_thermal_starting_index'
  'param_net1_buy_electric_starting_index'
  'param_net1_sell_electric_starting_index'
  'param_net1_emissions_starting_index'
  'param_P_net1_bat2_starting_index'
  'param_P_net1_heat_pump2_starting_index'
  'param_P_net1_bat1_starting_index'
  'param_P_net1_charging_station1_starting_index'
  'param_P_to_net1_starting_index'
  'param_P_net1_demand2_starting_index'
  'param_P_from_net1_starting_index'
  'param_P_net1_demand1_starting_index'
  'param_Q_from_net1_starting_index'
  'param_Q_to_net1_starting_index'
  'param_net1_inv_cost_starting_index'
  'param_P_net1_charging_station2_starting_index'
  'param_P_net1_heat_pump1_starting_index'
  'param_Q_net1_demand1_starting_index'
  'param_net1_buy_thermal_starting_index'
  'param_net2_buy_electric_starting_index'
  'param_P_net2_bat1_starting_index'
  'param_net2_inv_cost_starting_index'
  'param_net2_buy_thermal_starting_index'
  'param_net2_sell_thermal_starting_index'
  'param_P_net2_bat2_starting_index'
  'param_P_net2_charging_station1_starting_index'
  'param_Q_from_net2_starting_index'
  'param_P_net2_heat_pump2_starting_index'
  'param_P_from_net2_starting_index'
  'param_P_net2_demand2_starting_index'
  'param_Q_net2_demand2_starting_index'
  'param_net2_sell_electric_starting_index'
  'param_P_to_net2_starting_index'
  'param_P_net2_charging_station2_starting_index'
  'param_Q_to_net2_starting_index'
  'param_net2_emissions_starting_index'
  'param_P_net2_heat_pump1_starting_index'
  'param_P_net2_demand1_starting_index'
  'param_Q_net2_demand1_starting_index'
  'param_P_pv1_demand2_starting_index'
  'param_pv1_op_cost_starting_index'
  'param_pv1_emissions_starting_index'
  'param_P_pv1_demand1_starting_index'
  'param_P_pv1_bat1_starting_index'
  'param_P_pv1_charging_station2_starting_index'
  'param_P_pv1_bat2_starting_index'
  'param_P_pv1_net2_starting_index'
  'param_P_pv1_charging_station1_starting_index'
  'param_P_from_pv1_starting_index'
  'param_P_pv1_heat_pump2_starting_index'
  'param_P_pv1_heat_pump1_starting_index'
  'param_P_pv1_net1_starting_index'
  'param_pv1_inv_cost_starting_index'
  'param_P_pv2_charging_station1_starting_index'
  'param_P_pv2_net1_starting_index'
  'param_P_from_pv2_starting_index'
  'param_P_pv2_bat1_starting_index'
  'param_P_pv2_bat2_starting_index'
  'param_P_pv2_demand2_starting_index'
  'param_P_pv2_heat_pump1_starting_index'
  'param_P_pv2_heat_pump2_starting_index'
  'param_pv2_inv_cost_starting_index'
  'param_P_pv2_charging_station2_starting_index'
  'param_P_pv2_demand1_starting_index'
  'param_P_pv2_net2_starting_index'
  'param_pv2_op_cost_starting_index'
  'param_pv2_emissions_starting_index'
  'param_bat1_K_ch_starting_index'
  'param_P_bat1_net1_starting_index'
  'param_P_bat1_heat_pump1_starting_index'
  'param_P_bat1_demand2_starting_index'
  'param_bat1_cumulated_aging_starting_index'
  'param_bat1_op_cost_starting_index'
  'param_bat1_emissions_starting_index'
  'param_bat1_K_dis_starting_index'
  'param_bat1_integer_starting_index'
  'param_bat1_SOC_starting_index'
  'param_P_bat1_charging_station1_starting_index'
  'param_P_bat1_heat_pump2_starting_index'
  'param_P_to_bat1_starting_index'
  'param_bat1_inv_cost_starting_index'
  'param_P_from_bat1_starting_index'
  'param_P_bat1_demand1_starting_index'
  'param_P_bat1_charging_station2_starting_index'
  'param_P_bat1_net2_starting_index'
  'param_bat1_SOC_max_starting_index'
  'param_P_bat2_net2_starting_index'
  'param_P_bat2_demand1_starting_index'
  'param_bat2_SOC_starting_index'
  'param_P_bat2_charging_station1_starting_index'
  'param_bat2_op_cost_starting_index'
  'param_P_bat2_net1_starting_index'
  'param_bat2_inv_cost_starting_index'
  'param_bat2_emissions_starting_index'
  'param_bat2_cumulated_aging_starting_index'
  'param_P_bat2_demand2_starting_index'
  'param_P_to_bat2_starting_index'
  'param_P_bat2_charging_station2_starting_index'
  'param_P_bat2_heat_pump1_starting_index'
  'param_P_bat2_heat_pump2_starting_index'
  'param_P_from_bat2_starting_index'
  'param_bat2_K_dis_starting_index'
  'param_bat2_K_ch_starting_index'
  'param_bat2_SOC_max_starting_index'
  'param_bat2_integer_starting_index'
  'param_Q_CHP1_demand1_starting_index'
  'param_P_CHP1_bat2_starting_index'
  'param_P_CHP1_charging_station2_starting_index'
  'param_P_CHP1_demand2_starting_index'
  'param_CHP1_emissions_starting_index'
  'param_P_CHP1_heat_pump1_starting_index'
  'param_P_CHP1_bat1_starting_index'
  'param_P_CHP1_net2_starting_index'
  'param_CHP1_inv_cost_starting_index'
  'param_P_CHP1_net1_starting_index'
  'param_Q_CHP1_net2_starting_index'
  'param_P_CHP1_charging_station1_starting_index'
  'param_P_CHP1_heat_pump2_starting_index'
  'param_P_from_CHP1_starting_index'
  'param_Q_CHP1_net1_starting_index'
  'param_P_CHP1_demand1_starting_index'
  'param_CHP1_fuel_cons_starting_index'
  'param_CHP1_op_cost_starting_index'
  'param_Q_from_CHP1_starting_index'
  'param_Q_CHP1_demand2_starting_index'
  'param_P_CHP2_net2_starting_index'
  'param_P_CHP2_bat2_starting_index'
  'param_P_from_CHP2_starting_index'
  'param_CHP2_op_cost_starting_index'
  'param_Q_CHP2_net1_starting_index'
  'param_P_CHP2_heat_pump2_starting_index'
  'param_Q_CHP2_net2_starting_index'
  'param_Q_CHP2_demand1_starting_index'
  'param_CHP2_fuel_cons_starting_index'
  'param_P_CHP2_charging_station1_starting_index'
  'param_CHP2_inv_cost_starting_index'
  'param_P_CHP2_bat1_starting_index'
  'param_P_CHP2_demand1_starting_index'
  'param_Q_from_CHP2_starting_index'
  'param_Q_CHP2_demand2_starting_index'
  'param_P_CHP2_heat_pump1_starting_index'
  'param_P_CHP2_net1_starting_index'
  'param_CHP2_emissions_starting_index'
  'param_P_CHP2_demand2_starting_index'
  'param_P_CHP2_charging_station2_starting_index'
  'param_Q_solar_th1_net2_starting_index'
  'param_Q_from_solar_th1_starting_index'
  'param_Q_solar_th1_demand1_starting_index'
  'param_Q_solar_th1_net1_starting_index'
  'param_solar_th1_op_cost_starting_index'
  'param_solar_th1_inv_cost_starting_index'
  'param_solar_th1_emissions_starting_index'
  'param_Q_solar_th1_demand2_starting_index'
  'param_Q_solar_th2_net1_starting_index'
  'param_Q_solar_th2_net2_starting_index'
  'param_Q_solar_th2_demand2_starting_index'
  'param_solar_th2_emissions_starting_index'
  'param_Q_from_solar_th2_starting_index'
  'param_solar_th2_op_cost_starting_index'
  'param_Q_solar_th2_demand1_starting_index'
  'param_solar_th2_inv_cost_starting_index'
  'param_P_pvt1_bat2_starting_index'
  'param_P_from_pvt1_starting_index'
  'param_pvt1_emissions_starting_index'
  'param_pvt1_inv_cost_starting_index'
  'param_P_pvt1_bat1_starting_index'
  'param_P_pvt1_net2_starting_index'
  'param_pvt1_op_cost_starting_index'
  'param_P_pvt1_heat_pump2_starting_index'
  'param_P_pvt1_charging_station2_starting_index'
  'param_Q_pvt1_net1_starting_index'
  'param_Q_from_pvt1_starting_index'
  'param_P_pvt1_net1_starting_index'
  'param_Q_pvt1_demand2_starting_index'
  'param_P_pvt1_demand1_starting_index'
  'param_P_pvt1_heat_pump1_starting_index'
  'param_Q_pvt1_net2_starting_index'
  'param_P_pvt1_charging_station1_starting_index'
  'param_Q_pvt1_demand1_starting_index'
  'param_P_pvt1_demand2_starting_index'
  'param_P_pvt2_bat1_starting_index'
  'param_P_pvt2_demand2_starting_index'
  'param_P_pvt2_net1_starting_index'
  'param_Q_from_pvt2_starting_index'
  'param_P_from_pvt2_starting_index'
  'param_P_pvt2_charging_station2_starting_index'
  'param_P_pvt2_heat_pump2_starting_index'
  'param_P_pvt2_bat2_starting_index'
  'param_pvt2_emissions_starting_index'
  'param_Q_pvt2_demand1_starting_index'
  'param_P_pvt2_charging_station1_starting_index'
  'param_P_pvt2_demand1_starting_index'
  'param_Q_pvt2_net2_starting_index'
  'param_P_pvt2_net2_starting_index'
  'param_P_pvt2_heat_pump1_starting_index'
  'param_Q_pvt2_net1_starting_index'
  'param_pvt2_inv_cost_starting_index'
  'param_Q_pvt2_demand2_starting_index'
  'param_pvt2_op_cost_starting_index'
  'param_charging_station1_inv_cost_starting_index'
  'param_charging_station1_op_cost_starting_index'
  'param_charging_station1_emissions_starting_index'
  'param_charging_station2_emissions_starting_index'
  'param_charging_station2_inv_cost_starting_index'
  'param_charging_station2_op_cost_starting_index'
  'param_Q_heat_pump1_net1_starting_index'
  'param_heat_pump1_op_cost_starting_index'
  'param_heat_pump1_emissions_starting_index'
  'param_Q_from_heat_pump1_starting_index'
  'param_Q_to_heat_pump1_starting_index'
  'param_P_from_heat_pump1_starting_index'
  'param_heat_pump1_inv_cost_starting_index'
  'param_Q_heat_pump1_net2_starting_index'
  'param_Q_heat_pump1_demand1_starting_index'
  'param_Q_heat_pump1_demand2_starting_index'
  'param_P_to_heat_pump1_starting_index'
  'param_P_to_heat_pump2_starting_index'
  'param_Q_from_heat_pump2_starting_index'
  'param_heat_pump2_emissions_starting_index'
  'param_P_from_heat_pump2_starting_index'
  'param_Q_heat_pump2_demand2_starting_index'
  'param_Q_heat_pump2_net1_starting_index'
  'param_Q_heat_pump2_net2_starting_index'
  'param_Q_heat_pump2_demand1_starting_index'
  'param_heat_pump2_op_cost_starting_index'
  'param_heat_pump2_inv_cost_starting_index'
  'param_Q_to_heat_pump2_starting_index'
  'param_total_emissions_starting_index'
  'param_total_sell_starting_index'
  'param_total_buy_starting_index'
  'param_total_operation_cost_starting_index'
)

$values = @(
  45
  0
  0
  0
  0
  0
  0
  423.4403522009333
  0
  657.6983044018665
  0
  0
  0
  0
  0
  500
  920.5225047846375
  420.5225047846375
  592.5237571769562
  0
  0
  0
  0
  592.5237571769562
  183.6823647248564
  0
  0
  0
  160
  0
  0
  0
  800
  0
  0
  0
  800
  0
  0
  0
  0
  192
  0
  0
  0
  0
  1
  0
  0
  0
  0
  0
  0
  0
  0
  0
  0
  0
  0
  0
  0
  1.8
  0
  0
  0
  0.8999999999999999
  0.8999999999999999
  0
  0
  0
  0
  1
  0.8999999999999999
  1
  0
  0
  0
  0.000003166666666666667
  1
  0
  -0
  -0
  0.3
  0
  0
  0
  0
  0
  -0
  0
  0
  0.9999968333333333
  0
  0
  0.5
  0
  1
  0
  0
  0
  0
  -0
  0
  0
  0
  0
  0
  -0
  1
  1
  -0
  40
  0
  0
  0
  4.83
  9.25
  0
  0
  0
  0
  0
  0
  9.25
  20
  0
  1.5
  2.1
  10.5
  40
  0
  0
  0
  20
  10.5
  0
  9.25
  0
  0
  2.1
  0
  0
  0
  1.5
  40
  40
  9.25
  0
  4.83
  0
  0
  0
  0
  0
  0
  1
  0
  0
  0
  0
  0
  0
  0.6000000000000001
  1.2
  1
  1.2
  0
  0
  0
  0
  0
  0
  0
  1
  0
  0
  0
  0
  0
  0
  0
  0
  0
  0
  0
  0
  0
  0
  0
  1.56
  1.2
  0
  0.6000000000000001
  0
  0.7800000000000001
  1.56
  0
  0
  0
  0
  0.6000000000000001
  0
  0
  0
  1
  0
  0
  0
  0
  0
  0
  0
  8.561643835616438
  2.76
  80
  0
  0
  0
  0
  0
  80
  20
  20
  80
  2.76
  0
  80
  0
  0
  0
  8.561643835616438
  0
  0
  865.7783044018665
  0
  767.1227169257897
  646.7137533492462
)

for ($i = 0; $i -lt $names.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $names[$i]
  $ws.Cells.Item($row, 2).Value = $values[$i]
}

Write-Output "done"